# Update countries & provincias Spain
# This script applies a refreshed COVID-19 data snapshot to the "Pais"
# worksheet: several countries changed rank (causing their row to now
# display a different country name, since ranks/rows are fixed but the
# country occupying that rank changed) and numeric stats were refreshed
# for a handful of rows. The timestamp caption is also updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title / timestamp caption (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 15 de Junio de 2020 a las 07:47"

# --- Row 4: Estados Unidos (stats refresh, no rank change) ---
$ws.Range("B4").Value = 2162228
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 870050
$ws.Range("E4").Value = 1174320
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 117858

# --- Row 5: Brasil (stats refresh, no rank change) ---
$ws.Range("B5").Value = 867882
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 453568
$ws.Range("E5").Value = 370925
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 43389

# --- Rows 65-68: Camerun overtakes Guatemala, Honduras and Marruecos ---
$ws.Range("A65").Value = "Camerun"
$ws.Range("B65").Value = 9864
$ws.Range("C65").Value = 1183
$ws.Range("D65").Value = 5570
$ws.Range("E65").Value = 4018
$ws.Range("F65").Value = 0
$ws.Range("G65").Value = 64
$ws.Range("H65").Value = 276

$ws.Range("A66").Value = "Guatemala"
$ws.Range("B66").Value = 9845
$ws.Range("C66").Value = 354
$ws.Range("D66").Value = 1886
$ws.Range("E66").Value = 7575
$ws.Range("F66").Value = 0
$ws.Range("G66").Value = 17
$ws.Range("H66").Value = 384

$ws.Range("A67").Value = "Honduras"
$ws.Range("B67").Value = 8858
$ws.Range("C67").Value = 403
$ws.Range("D67").Value = 967
$ws.Range("E67").Value = 7579
$ws.Range("F67").Value = 0
$ws.Range("G67").Value = 2
$ws.Range("H67").Value = 312

$ws.Range("A68").Value = "Marruecos"
$ws.Range("B68").Value = 8793
$ws.Range("C68").Value = 0
$ws.Range("D68").Value = 7765
$ws.Range("E68").Value = 816
$ws.Range("F68").Value = 0
$ws.Range("G68").Value = 0
$ws.Range("H68").Value = 212

# --- Rows 75-77: Uzbekistan overtakes Senegal and Costa de Marfil ---
$ws.Range("A75").Value = "Uzbekistan"
$ws.Range("B75").Value = 5103
$ws.Range("C75").Value = 23
$ws.Range("D75").Value = 3943
$ws.Range("E75").Value = 1141
$ws.Range("F75").Value = 0
$ws.Range("G75").Value = 0
$ws.Range("H75").Value = 19

$ws.Range("A76").Value = "Senegal"
$ws.Range("B76").Value = 5090
$ws.Range("C76").Value = 0
$ws.Range("D76").Value = 3344
$ws.Range("E76").Value = 1686
$ws.Range("F76").Value = 0
$ws.Range("G76").Value = 0
$ws.Range("H76").Value = 60

$ws.Range("A77").Value = "Costa de Marfil"
$ws.Range("B77").Value = 5084
$ws.Range("C77").Value = 0
$ws.Range("D77").Value = 2505
$ws.Range("E77").Value = 2534
$ws.Range("F77").Value = 0
$ws.Range("G77").Value = 0
$ws.Range("H77").Value = 45

# --- Row 146: Togo (stats refresh, no rank change) ---
$ws.Range("B146").Value = 531
$ws.Range("C146").Value = 1
$ws.Range("D146").Value = 299
$ws.Range("E146").Value = 219

# --- Rows 206-209: Groenlandia overtakes Islas Malvinas;
#     Islas Turcas y Caicos overtakes Santa Sede ---
$ws.Range("A206").Value = "Groenlandia"
$ws.Range("A207").Value = "Islas Malvinas"

$ws.Range("A208").Value = "Islas Turcas y Caicos"
$ws.Range("B208").Value = 12
$ws.Range("C208").Value = 0
$ws.Range("D208").Value = 11
$ws.Range("E208").Value = 0
$ws.Range("F208").Value = 0
$ws.Range("G208").Value = 0
$ws.Range("H208").Value = 1

$ws.Range("A209").Value = "Santa Sede"
$ws.Range("B209").Value = 12
$ws.Range("C209").Value = 0
$ws.Range("D209").Value = 12
$ws.Range("E209").Value = 0
$ws.Range("F209").Value = 0
$ws.Range("G209").Value = 0
$ws.Range("H209").Value = 0

$wb.Save()
